# Auto-generated edit script: updates crypto price/volume table
# to reflect refreshed data from GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be written as literal TEXT even when the
# string looks like a number (e.g. "58.53"), then restore the
# cell style back to the default "Normal" style so no extra
# number-format/style is introduced.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "36.363.13"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "1.963.64"
$ws.Range("E3").Value = "  -4.87%  "
Set-TextCell $ws.Range("D5") "243.93"
$ws.Range("E5").Value = "  -2.89%  "
Set-TextCell $ws.Range("D6") "0.619"
$ws.Range("E6").Value = "  -4.41%  "
Set-TextCell $ws.Range("D7") "58.53"
$ws.Range("E7").Value = "  -9.33%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextCell $ws.Range("D9") "0.372"
$ws.Range("E9").Value = "  -5.49%  "
Set-TextCell $ws.Range("D10") "55.51"
$ws.Range("E10").Value = "  -6.50%  "
Set-TextCell $ws.Range("D11") "0.0848"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  -0.73%  "
Set-TextCell $ws.Range("D13") "22.03"
$ws.Range("E13").Value = "  -4.61%  "
Set-TextCell $ws.Range("D14") "0.829"
$ws.Range("E14").Value = "  -9.16%  "
$ws.Range("D15").Value = "2.251.98"
$ws.Range("E15").Value = "  -4.64%  "
Set-TextCell $ws.Range("D16") "13.48"
$ws.Range("E16").Value = "  -8.03%  "
Set-TextCell $ws.Range("D17") "5.33"
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("D18").Value = "1.976.31"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "36.277.18"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D21") "70.63"
$ws.Range("E21").Value = "  -3.70%  "
Set-TextCell $ws.Range("D22") "230.61"
$ws.Range("E22").Value = "  -3.31%  "
Set-TextCell $ws.Range("D23") "5.10"
$ws.Range("E23").Value = "  -6.74%  "
$ws.Range("E24").Value = "  -0.09%  "
Set-TextCell $ws.Range("D25") "2.50"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  -4.07%  "
Set-TextCell $ws.Range("D27") "9.59"
$ws.Range("E27").Value = "  -3.37%  "
Set-TextCell $ws.Range("D28") "166.37"
$ws.Range("E28").Value = "  +3.49%  "
Set-TextCell $ws.Range("D29") "19.72"
$ws.Range("E29").Value = "  -4.24%  "
Set-TextCell $ws.Range("D30") "0.117"
$ws.Range("E30").Value = "  -12.59%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("E32").Value = "  -3.09%  "
Set-TextCell $ws.Range("D33") "4.74"
$ws.Range("E33").Value = "  -7.55%  "
Set-TextCell $ws.Range("D34") "0.0637"
$ws.Range("E34").Value = "  +1.96%  "
Set-TextCell $ws.Range("D35") "4.31"
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws.Range("D37") "6.05"
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D38") "1.81"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -12.03%  "
Set-TextCell $ws.Range("D40") "2.87"
$ws.Range("E40").Value = "  -4.64%  "
Set-TextCell $ws.Range("D41") "0.0970"
$ws.Range("E41").Value = "  -4.77%  "
$ws.Range("E42").Value = "  -4.68%  "
Set-TextCell $ws.Range("D43") "1.18"
$ws.Range("E43").Value = "  -7.56%  "
Set-TextCell $ws.Range("D44") "0.0210"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws.Range("D45") "15.80"
$ws.Range("E45").Value = "  -7.88%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Range("D46") "1.05"
$ws.Range("E46").Value = "  -9.78%  "
Set-TextCell $ws.Range("D47") "88.75"
$ws.Range("E47").Value = "  -6.54%  "
$ws.Range("D48").Value = "1.344.07"
$ws.Range("E48").Value = "  -3.24%  "
Set-TextCell $ws.Range("D49") "7.25"
$ws.Range("E49").Value = "  -7.55%  "
Set-TextCell $ws.Range("D50") "2.82"
$ws.Range("E50").Value = "  -3.55%  "
Set-TextCell $ws.Range("D51") "44.70"
$ws.Range("E51").Value = "  -3.05%  "
